# Applies the "created only bootstrap with more candidates and examples" edit:
#  - Adds two new header columns AX1/AY1: "Max Bootstrapped Demos" and
#    "Number of Candidate Programs" (matching the bold/bordered header style
#    used by the rest of row 1).
#  - Appends a new data row (row 12) for the qwen2:7b-instruct-q5_K_M /
#    llama3:70b run, including values for the two new columns.
#  - The sheet's used range grows from A1:AW11 to A1:AY12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (AX1, AY1), with the same formatting as the existing
#     bold/bordered/centered header cells (copy format from AW1). ---
$ws.Range("AW1").Copy() | Out-Null
$ws.Range("AX1:AY1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("AX1").Value2 = "Max Bootstrapped Demos"
$ws.Range("AY1").Value2 = "Number of Candidate Programs"

# --- New data row 12 ---
$ws.Range("A12").Value2 = "qwen2:7b-instruct-q5_K_M"
$ws.Range("B12").Value2 = "llama3:70b"
$ws.Range("C12").Value2 = 1
$ws.Range("D12").Value2 = 200
$ws.Range("E12").Value2 = 2138.53
$ws.Range("F12").Value2 = 50.6
$ws.Range("G12").Value2 = 42.5
$ws.Range("H12").Value2 = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_val_match.txt"
$ws.Range("I12").Value2 = 50.6
$ws.Range("J12").Value2 = 90
$ws.Range("K12").Value2 = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_val_correct.txt"
$ws.Range("L12").Value2 = 31.33333333333333
$ws.Range("M12").Value2 = 94.84999999999999
$ws.Range("N12").Value2 = 51.25
$ws.Range("O12").Value2 = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_test_match.txt"
$ws.Range("P12").Value2 = 94.84999999999999
$ws.Range("Q12").Value2 = 91.25
$ws.Range("R12").Value2 = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_test_correct.txt"
$ws.Range("S12").Value2 = 42.33333333333334
$ws.Range("T12").Value2 = 0
$ws.Range("U12").Value2 = 92.77
$ws.Range("V12").Value2 = 57.5
$ws.Range("W12").Value2 = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_val_fewshot_match.txt"
$ws.Range("X12").Value2 = 92.77
$ws.Range("Y12").Value2 = 90
$ws.Range("Z12").Value2 = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_val_fewshot_correct.txt"
$ws.Range("AA12").Value2 = 40.66666666666666
$ws.Range("AB12").Value2 = 179.56
$ws.Range("AC12").Value2 = 43.75
$ws.Range("AD12").Value2 = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_test_fewshot_match.txt"
$ws.Range("AE12").Value2 = 179.56
$ws.Range("AF12").Value2 = 88.75
$ws.Range("AG12").Value2 = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_test_fewshot_correct.txt"
$ws.Range("AH12").Value2 = 31.33333333333333
$ws.Range("AI12").Value2 = 1138.61
$ws.Range("AJ12").Value2 = 199.39
$ws.Range("AK12").Value2 = 47.5
$ws.Range("AL12").Value2 = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_val_bootstrap_match.txt"
$ws.Range("AM12").Value2 = 199.39
$ws.Range("AN12").Value2 = 82.5
$ws.Range("AO12").Value2 = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_val_bootstrap_correct.txt"
$ws.Range("AP12").Value2 = 31.33333333333333
$ws.Range("AQ12").Value2 = 382.76
$ws.Range("AR12").Value2 = 46.25
$ws.Range("AS12").Value2 = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_test_bootstrap_match.txt"
$ws.Range("AT12").Value2 = 382.76
$ws.Range("AU12").Value2 = 85
$ws.Range("AV12").Value2 = "logs\qwen2_7b_instruct_q5_K_M_llama3_70b_1_200_test_bootstrap_correct.txt"
$ws.Range("AW12").Value2 = 31
$ws.Range("AX12").Value2 = 2
$ws.Range("AY12").Value2 = 2

Write-Output "Applied log_evaluations edit: added AX/AY headers + row 12"
